# Rebuild Sheet1 content according to the target layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe all existing cell contents (old headers / FUNCEXEC demo rows).
$ws.Cells.ClearContents()

# Row 1 - new headers
$ws.Range("A1").Value = "c(a)"
$ws.Range("B1").Value = "a"

# Row 2 - formula referencing B2, and input value
$ws.Range("A2").Formula = '=FUNCEXEC("c_from_a",B2)'
$ws.Range("B2").Value = 1

# Row 3 - extra input value
$ws.Range("B3").Value = 5

# Row 4 - extra input value
$ws.Range("B4").Value = 10

# Row 7 - formula referencing B10 (outside populated data)
$ws.Range("A7").Formula = '=FUNCEXEC("c_from_a",B10)'

# Selection ends up on A8, as in the target workbook
$ws.Range("A8").Select()

$wb.Save()
